$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the TC_ID values in column B from SCD0338-0XX to SCD0026-0XX
$ws.Range("B2").Value = "SCD0026-003"
$ws.Range("B3").Value = "SCD0026-004"
$ws.Range("B4").Value = "SCD0026-005"
$ws.Range("B5").Value = "SCD0026-006"
$ws.Range("B6").Value = "SCD0026-007"

# Update the active selection to B7
$ws.Range("B7").Select()
